$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")
$ws.Activate()

# New work-log entry (row 20) on the "Eetu Pihamäki" sheet.
$ws.Range("A20").Value = 43395
$ws.Range("B20").Value = 0.71875
$ws.Range("C20").Value = 0.80555555555555547
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = "5 min esityslistan teko ja lähetys. 1h 20 min AD Windows Server 2016  pavelimen asennusta ja konfigurointia VirtualBox ympäristössä. 30 min valmistautumista HTTPS-yhteyden konfigurointiin midPoint IdM-palvelimeen (SSL, encryption keys, Java JCE Keystore jne.) https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2022.10.2018.txt"

# Row grew taller to fit the wrapped text of the new note.
$ws.Rows(20).RowHeight = 120

# Scroll / selection state left by the author after entering the row.
$ws.Range("G20").Select()
